$wb = $excel.ActiveWorkbook

# --- Rename existing sheet, add the new "metadatos" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws1.Name = "datos"
$ws2.Name = "metadatos"

# --- Give the new sheet's cells the explicit (non-theme) Calibri font used
#     throughout the metadata table, matching the style the source workbook
#     ships with (fontId 1, no scheme). ---
$ws2.Range("A1:D7").Font.Name = "Calibri"

# --- Header row ---
$ws2.Range("A1").Value = "Variables"
$ws2.Range("B1").Value = "Descripción"
$ws2.Range("C1").Value = "Fuente"
$ws2.Range("D1").Value = "Fecha_de_extracción"

# --- Apply the date number format to the extraction-date column before
#     writing the values so the engine maps straight to the builtin
#     numFmtId 15 (d-mmm-yy) instead of synthesizing a custom format. ---
$ws2.Range("D2:D6").NumberFormat = "d-mmm-yy"

# --- Data rows ---
$ws2.Range("A2").Value = "anno"
$ws2.Range("B2").Value = "Año"
$ws2.Range("C2").Value = "…"
$ws2.Range("D2").Value = "2025-03-06"

$ws2.Range("A3").Value = "codmpio"
$ws2.Range("B3").Value = "Código del municipio"
$ws2.Range("C3").Value = "…"
$ws2.Range("D3").Value = "2025-03-06"

$ws2.Range("A4").Value = "numerador"
$ws2.Range("B4").Value = "# de casos de Violencia Intrafamiliar en niños, niñas y adolescentes"
$ws2.Range("C4").Value = "Instituto Nacional de Medicina Legal y Ciencias Forenses"
$ws2.Range("D4").Value = "2025-03-06"

$ws2.Range("A5").Value = "denominador"
$ws2.Range("B5").Value = "Total niños/niñas/adolescentes"
$ws2.Range("C5").Value = "Departamento Administrativo Nacional de Estadística (DANE)"
$ws2.Range("D5").Value = "2025-03-06"

$ws2.Range("A6").Value = "intrafamiliar"
$ws2.Range("B6").Value = "x 100,000"
$ws2.Range("C6").Value = "Elaboración Propia"
$ws2.Range("D6").Value = "2025-03-06"

# Row 7 stays empty, but still carries the sheet's default style.
# (Range.Font.Name above already stamped the style onto these cells.)

# --- Selections: datos -> C3, metadatos -> D2:D6 (metadatos ends up the
#     active/visible tab, matching tabSelected moving sheets). ---
[void]$ws1.Range("C3").Select()
[void]$ws2.Range("D2:D6").Select()
